$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Resize the table to include the new column H first (creates default "Column8" header)
$tbl = $ws.ListObjects.Item("Tabla1")
$tbl.Resize($ws.Range("A1:H9"))

# Now set header for new column H - this also renames the table column
$ws.Range("H1").Value = "total_clp"

# Update column G values (divide previous hours-like measure into turnos count) and add H totals
$ws.Range("G2").Value = 6
$ws.Range("H2").Value = 15962400

$ws.Range("G3").Value = 6
$ws.Range("H3").Value = 14366160

$ws.Range("G4").Value = 6
$ws.Range("H4").Value = 11447849

$ws.Range("G5").Value = 6
$ws.Range("H5").Value = 5724294

$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0

$ws.Range("G7").Value = 6
$ws.Range("H7").Value = 10567700

$ws.Range("G8").Value = 6
$ws.Range("H8").Value = 5283850

$ws.Range("G9").Value = 6
$ws.Range("H9").Value = 5283850

# Set column H width to match target (~11.27 chars in OOXML units)
$ws.Columns.Item(8).ColumnWidth = 10.43

# Set selection to match final state
$ws.Range("L6").Select()

$wb.Save()
